$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the data rows as a cross-product of source wells (A1/A2/A3,
# each with its own transfer volume + reagent) against destination
# wells (A1/A2/A3), replacing the original 3-row table with 9 rows.

$sourcePlateName = "level 2 LDV source plate"
$sourcePlateType = "384LDV_AQ_B"
$destPlateName = "384-Well Level 2 MoClo output plate"
$destPlateType = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"

$sources = @(
    @{ Well = "A1"; Volume = 500; Reagent = "DNA ligase buffer" },
    @{ Well = "A2"; Volume = 125; Reagent = "DNA ligase" },
    @{ Well = "A3"; Volume = 250; Reagent = "BsmBI (NEB)" }
)
$destWells = @("A1", "A2", "A3")

$row = 2
$uid = 1
foreach ($src in $sources) {
    foreach ($dw in $destWells) {
        $ws.Cells.Item($row, 1).Value = $uid
        $ws.Cells.Item($row, 2).Value = $sourcePlateName
        $ws.Cells.Item($row, 3).Value = $sourcePlateType
        $ws.Cells.Item($row, 4).Value = $src.Well
        $ws.Cells.Item($row, 5).Value = $destPlateName
        $ws.Cells.Item($row, 6).Value = $destPlateType
        $ws.Cells.Item($row, 7).Value = $dw
        $ws.Cells.Item($row, 8).Value = $src.Volume
        $ws.Cells.Item($row, 9).Value = $src.Reagent

        $uid++
        $row++
    }
}
